{"js": "// Remove the 404-screenshot paragraph (the inline picture right under the\n// \"Introducci\u00f3n\" heading) together with the descriptive paragraph that\n// follows it (\"Este proyecto consiste en la creaci\u00f3n de un motor de\n// plantillas tipo Pug escrito \u00edntegramente en PHP...\"). The heading itself\n// and the paragraph that used to come after the description (\"El proyecto\n// evoluciona...\") are left untouched and end up adjacent to the heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Introducci\u00f3n\" heading paragraph.\nlet introIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Introducci\u00f3n\") {\n    introIndex = i;\n    break;\n  }\n}\n\nif (introIndex === -1) {\n  throw new Error('Paragraph \"Introducci\u00f3n\" not found.');\n}\n\n// Load the inline-picture collections of the couple of paragraphs that\n// immediately follow the heading so we can positively identify the\n// screenshot paragraph instead of guessing from its (empty) text.\nconst candidates = items.slice(introIndex + 1, introIndex + 3);\ncandidates.forEach((p) => p.inlinePictures.load(\"items\"));\nawait context.sync();\n\nconst toDelete = [];\nfor (const p of candidates) {\n  const text = p.text.trim();\n  const hasPicture = p.inlinePictures.items.length > 0;\n  if (hasPicture || text.startsWith(\"Este proyecto consiste en la creaci\u00f3n de un\")) {\n    toDelete.push(p);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the 404-screenshot paragraph (inline picture) and the descriptive\n# paragraph right after it (\"Este proyecto consiste en la creaci\u00f3n de un\n# motor de plantillas tipo Pug ...\"), both of which sit directly under the\n# \"Introducci\u00f3n\" heading. The heading itself and the paragraph that follows\n# (\"El proyecto evoluciona...\") are left untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Introducci\u00f3n\" heading paragraph (1-based COM index).\n$introIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Introducci\u00f3n\") {\n        $introIndex = $i\n    }\n}\n\nif ($introIndex -eq -1) {\n    throw \"Paragraph 'Introduccion' not found.\"\n}\n\n$imgIndex = $introIndex + 1\n$textIndex = $introIndex + 2\n\n# Delete any inline picture anchored inside the image paragraph. (Range.InlineShapes\n# can mis-resolve, so look the shape up through the document-wide collection and\n# match it by start position instead of indexing into the range sub-collection.)\n$imgPara = $d.Paragraphs.Item($imgIndex)\n$pStart = $imgPara.Range.Start\n$pEnd = $imgPara.Range.End\n\n$targetShape = $null\nforeach ($shp in $d.InlineShapes) {\n    $shpStart = $shp.Range.Start\n    if ($shpStart -ge $pStart -and $shpStart -lt $pEnd) {\n        $targetShape = $shp\n    }\n}\nif ($targetShape -ne $null) {\n    $targetShape.Delete()\n}\n\n# Delete the two paragraphs (now-empty image paragraph + the descriptive\n# paragraph) in reverse order so earlier indices stay valid.\n$d.Paragraphs.Item($textIndex).Range.Delete()\n$d.Paragraphs.Item($imgIndex).Range.Delete()\n"}
